# Re-arranged contrasts in forest plot
# The data table (A1:T17, header in row 1) is sorted by column C (year), ascending,
# using a stable sort so that rows with identical years keep their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the worksheet Sort object (rather than Range.Sort) so that Excel records the
# sort state (sortState/sortCondition) in the saved worksheet, matching a manual
# Data > Sort operation performed by a user.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("C1:C17")) | Out-Null
$sortObj.SetRange($ws.Range("A1:T17"))
$sortObj.Header = 1
$sortObj.MatchCase = $false
$sortObj.Apply()

# Reflect the final view state left behind after the sort: the window is scrolled
# back to show column A (no frozen/scrolled topLeftCell) and the active selection
# is on cell F12.
$ws.Range("A1").Select() | Out-Null
$ws.Range("F12").Select() | Out-Null
